$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G9 needs to hold the literal text "0.01" (same shared string already used
# by J2), without altering the cell's existing style (s="4", General format).
# A direct Value/Formula assignment would get re-parsed as a number (or pick
# up a new Text-formatted style via NumberFormat), so copy the existing text
# cell J2 and paste only its value into G9 - this matches how Excel keeps
# the cell's style untouched while making the content a shared string.
$ws.Range("J2").Copy()
$ws.Range("G9").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

$ws.Range("G10").Select()
